$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 32.9414084707114
$ws.Range("B1").Value = 108.116717780063
$ws.Range("C1").Value = 188.851263741428
$ws.Range("D1").Value = 9.29911081180866
$ws.Range("E1").Value = 91.1148366942605
$ws.Range("F1").Value = 156.880231554099
$ws.Range("G1").Value = 11.7141144404719
$ws.Range("H1").Value = 149.384023877505
$ws.Range("I1").Value = 53.415918561358
$ws.Range("J1").Value = 164.375875314873
$ws.Range("A2").Value = 28.8470517046969
$ws.Range("B2").Value = 140.439572716337
$ws.Range("C2").Value = 124.278133792932
$ws.Range("D2").Value = 0.282906275374306
$ws.Range("E2").Value = 112.74384237488
$ws.Range("F2").Value = 86.9387512500113
$ws.Range("G2").Value = 150.847064308285
$ws.Range("H2").Value = 179.575945334312
$ws.Range("I2").Value = 107.160048516076
$ws.Range("J2").Value = 190.545255965807
$ws.Range("A3").Value = 66.217780423452
$ws.Range("B3").Value = 7.18497354871825
$ws.Range("C3").Value = 3.4049037859798
$ws.Range("D3").Value = 128.198529560211
$ws.Range("E3").Value = 174.958692851923
$ws.Range("F3").Value = 154.288465228997
$ws.Range("G3").Value = 25.108172942469
$ws.Range("H3").Value = 103.82220945499
$ws.Range("I3").Value = 132.524444410822
$ws.Range("J3").Value = 108.985220039722
$ws.Range("A4").Value = 195.135846172988
$ws.Range("B4").Value = 2.64957714949249
$ws.Range("C4").Value = 5.7952925589845
$ws.Range("D4").Value = 102.310986119467
$ws.Range("E4").Value = 18.1320874104891
$ws.Range("F4").Value = 114.406534989554
$ws.Range("G4").Value = 181.640655725096
$ws.Range("H4").Value = 7.13782022108222
$ws.Range("I4").Value = 152.050067648315
$ws.Range("J4").Value = 182.08702103332
$ws.Range("A5").Value = 183.363339576574
$ws.Range("B5").Value = 11.9944848176066
$ws.Range("C5").Value = 5.11998692765831
$ws.Range("D5").Value = 25.8805357971604
$ws.Range("E5").Value = 160.020234324047
$ws.Range("F5").Value = 131.233907645212
$ws.Range("G5").Value = 143.220596640939
$ws.Range("H5").Value = 89.440680336878
$ws.Range("I5").Value = 58.8428589789396
$ws.Range("J5").Value = 66.359319568779
$ws.Range("A6").Value = 172.180705783973
$ws.Range("B6").Value = 175.792729191386
$ws.Range("C6").Value = 27.6833604218826
$ws.Range("D6").Value = 24.3489646466211
$ws.Range("E6").Value = 84.2274857145862
$ws.Range("F6").Value = 25.7564349219931
$ws.Range("G6").Value = 104.711813994083
$ws.Range("H6").Value = 60.6527341812163
$ws.Range("I6").Value = 34.3404179598859
$ws.Range("J6").Value = 136.826371465263
$ws.Range("A7").Value = 131.77205861163
$ws.Range("B7").Value = 107.891904985482
$ws.Range("C7").Value = 16.8595794666836
$ws.Range("D7").Value = 144.430698521636
$ws.Range("E7").Value = 169.240029141884
$ws.Range("F7").Value = 26.1974745552044
$ws.Range("G7").Value = 134.644280157352
$ws.Range("H7").Value = 21.9671476734649
$ws.Range("I7").Value = 182.150818864885
$ws.Range("J7").Value = 198.337307385326
$ws.Range("A8").Value = 105.298095524916
$ws.Range("B8").Value = 143.709244087203
$ws.Range("C8").Value = 27.5258776859966
$ws.Range("D8").Value = 125.073027482756
$ws.Range("E8").Value = 7.18191638923339
$ws.Range("F8").Value = 54.2232956058454
$ws.Range("G8").Value = 2.06498662105994
$ws.Range("H8").Value = 177.524367988819
$ws.Range("I8").Value = 168.178295236164
$ws.Range("J8").Value = 43.7247852067113
$ws.Range("A9").Value = 11.0678685880582
$ws.Range("B9").Value = 135.667492605591
$ws.Range("C9").Value = 44.9793504760504
$ws.Range("D9").Value = 66.1651248420426
$ws.Range("E9").Value = 136.804514255749
$ws.Range("F9").Value = 19.3431169816028
$ws.Range("G9").Value = 174.96621672761
$ws.Range("H9").Value = 181.446327912363
$ws.Range("I9").Value = 18.0835004048811
$ws.Range("J9").Value = 192.375652488496
$ws.Range("A10").Value = 9.69472099547029
$ws.Range("B10").Value = 120.987921543879
$ws.Range("C10").Value = 172.797402261196
$ws.Range("D10").Value = 15.2236961830518
$ws.Range("E10").Value = 50.3149624216906
$ws.Range("F10").Value = 75.4714345910919
$ws.Range("G10").Value = 195.134905350923
$ws.Range("H10").Value = 60.6892884060225
$ws.Range("I10").Value = 56.6405066552761
$ws.Range("J10").Value = 133.822759768843
$ws.Range("A11").Value = 196.58962748786
$ws.Range("B11").Value = 121.253448967474
$ws.Range("C11").Value = 107.289861471993
$ws.Range("D11").Value = 60.5055515936136
$ws.Range("E11").Value = 161.061224043863
$ws.Range("F11").Value = 28.4714616967698
$ws.Range("G11").Value = 148.266851505389
$ws.Range("H11").Value = 102.610332939127
$ws.Range("I11").Value = 17.1670482573877
$ws.Range("J11").Value = 30.0041901087408
$ws.Range("A12").Value = 23.6914483009332
$ws.Range("B12").Value = 127.187446005264
$ws.Range("C12").Value = 92.4744389450524
$ws.Range("D12").Value = 190.615632753175
$ws.Range("E12").Value = 125.758502877205
$ws.Range("F12").Value = 196.104566006039
$ws.Range("G12").Value = 62.9125545094314
$ws.Range("H12").Value = 150.694454624641
$ws.Range("I12").Value = 7.62618426588652
$ws.Range("J12").Value = 149.896912160282
$ws.Range("A13").Value = 51.2312578275945
$ws.Range("B13").Value = 153.197952244989
$ws.Range("C13").Value = 3.88364726858383
$ws.Range("D13").Value = 189.775166376389
$ws.Range("E13").Value = 188.642586389856
$ws.Range("F13").Value = 184.310173981036
$ws.Range("G13").Value = 170.911841826007
$ws.Range("H13").Value = 12.3021815029449
$ws.Range("I13").Value = 74.7580650610654
$ws.Range("J13").Value = 131.710481798142
$ws.Range("A14").Value = 59.0883902549224
$ws.Range("B14").Value = 141.375698215037
$ws.Range("C14").Value = 120.883861333543
$ws.Range("D14").Value = 34.3555354673208
$ws.Range("E14").Value = 47.1351577188518
$ws.Range("F14").Value = 89.814419620584
$ws.Range("G14").Value = 28.3776311335981
$ws.Range("H14").Value = 184.473798882437
$ws.Range("I14").Value = 105.103900798179
$ws.Range("J14").Value = 108.33305255898
$ws.Range("A15").Value = 71.0762654762139
$ws.Range("B15").Value = 72.3558837884832
$ws.Range("C15").Value = 164.279279654976
$ws.Range("D15").Value = 188.07931029614
$ws.Range("E15").Value = 168.684204187563
$ws.Range("F15").Value = 82.5072749902062
$ws.Range("G15").Value = 28.513482598827
$ws.Range("H15").Value = 182.181769508022
$ws.Range("I15").Value = 89.4651933058469
$ws.Range("J15").Value = 54.210396415652
$ws.Range("A16").Value = 12.5588800816605
$ws.Range("B16").Value = 44.4404507262821
$ws.Range("C16").Value = 53.063104140136
$ws.Range("D16").Value = 106.743594494995
$ws.Range("E16").Value = 82.5915019412486
$ws.Range("F16").Value = 43.3916752428709
$ws.Range("G16").Value = 117.36980169889
$ws.Range("H16").Value = 117.514695095604
$ws.Range("I16").Value = 71.8629652037578
$ws.Range("J16").Value = 176.751050062827
$ws.Range("A17").Value = 57.5596198707631
$ws.Range("B17").Value = 135.964670002444
$ws.Range("C17").Value = 27.8522678780613
$ws.Range("D17").Value = 85.4565664592462
$ws.Range("E17").Value = 170.915799853818
$ws.Range("F17").Value = 82.3157500858958
$ws.Range("G17").Value = 6.30358467172067
$ws.Range("H17").Value = 58.1189034777316
$ws.Range("I17").Value = 143.480475034323
$ws.Range("J17").Value = 35.9440832566209
$ws.Range("A18").Value = 167.72693487244
$ws.Range("B18").Value = 78.4387556269945
$ws.Range("C18").Value = 45.5905538264618
$ws.Range("D18").Value = 99.293131706907
$ws.Range("E18").Value = 78.8206466840676
$ws.Range("F18").Value = 178.875374039111
$ws.Range("G18").Value = 188.918672590013
$ws.Range("H18").Value = 15.8043369724435
$ws.Range("I18").Value = 21.0909621888264
$ws.Range("J18").Value = 106.13531139965
$ws.Range("A19").Value = 155.796691382209
$ws.Range("B19").Value = 188.730072317985
$ws.Range("C19").Value = 122.836988197098
$ws.Range("D19").Value = 20.5476686454134
$ws.Range("E19").Value = 119.151601716481
$ws.Range("F19").Value = 14.6479395286403
$ws.Range("G19").Value = 88.3125940749015
$ws.Range("H19").Value = 14.1402668385488
$ws.Range("I19").Value = 151.764033526072
$ws.Range("J19").Value = 3.74348247598088
$ws.Range("A20").Value = 172.444617921694
$ws.Range("B20").Value = 110.862936037994
$ws.Range("C20").Value = 112.610833678679
$ws.Range("D20").Value = 128.352850735352
$ws.Range("E20").Value = 50.7734326882164
$ws.Range("F20").Value = 135.11159547377
$ws.Range("G20").Value = 44.5036159104219
$ws.Range("H20").Value = 78.8227131957294
$ws.Range("I20").Value = 17.1635104423219
$ws.Range("J20").Value = 182.283408093398
